$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 32 (which holds the 2021-12-10 entry).
# This pushes the existing row 32 -> 33 and row 33 -> 34, matching the diff.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new weekly data point.
$ws.Cells.Item(32, 1).Value = 11
$ws.Cells.Item(32, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(32, 3).Value = "Bíobío"
$ws.Cells.Item(32, 4).Value = Get-Date -Year 2022 -Month 1 -Day 7 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(32, 5).Value = 8
$ws.Cells.Item(32, 6).Value = 100112031
$ws.Cells.Item(32, 7).Value = "Poroto verde"
$ws.Cells.Item(32, 8).Value = "Magnum"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 200
$ws.Cells.Item(32, 11).Value = 25000
$ws.Cells.Item(32, 12).Value = 26000
$ws.Cells.Item(32, 13).Value = 25500
$ws.Cells.Item(32, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(32, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(32, 16).Value = 1020
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"
